$wb = $excel.ActiveWorkbook

# Sheet: VENTAS POR GRUPO
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsGrupo.Range("H26").Value = 853.2
$wsGrupo.Range("M45").Value = 2352.42
$wsGrupo.Range("H54").Value = "4 de 52"
$wsGrupo.Range("M54").Value = "15 de 52"

# Sheet: VENTA MENSUAL
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsMensual.Range("F26").Value = 3066.52
$wsMensual.Range("F45").Value = 3653.82
$wsMensual.Range("F58").Value = 50208.69

# Sheet: CUMPLIMIENTO MENSUAL
$wsCumpl = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$wsCumpl.Range("D6").Value = 3764.7
$wsCumpl.Range("E6").Value = -857.11631853974
$wsCumpl.Range("F6").Value = 1.294786466166049

$wsCumpl.Range("D12").Value = 26278.39
$wsCumpl.Range("E12").Value = 35585.3303947566
$wsCumpl.Range("F12").Value = 0.4247786882572824

$wsCumpl.Range("D15").Value = 49057.92999999999
$wsCumpl.Range("E15").Value = 72996.90551083436
$wsCumpl.Range("F15").Value = 0.4019335227045986
